# Auto-generated edit script: update "想去人数" (F column) counts
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 518
$ws.Range("F4").Value = 217
$ws.Range("F5").Value = 103
$ws.Range("F6").Value = 345
$ws.Range("F7").Value = 1324
$ws.Range("F10").Value = 1382
$ws.Range("F13").Value = 212
$ws.Range("F14").Value = 150
$ws.Range("F15").Value = 269
$ws.Range("F16").Value = 1718
$ws.Range("F18").Value = 288
$ws.Range("F19").Value = 339
$ws.Range("F20").Value = 3639
$ws.Range("F21").Value = 35
$ws.Range("F22").Value = 426
$ws.Range("F23").Value = 953
$ws.Range("F24").Value = 1231
$ws.Range("F26").Value = 2881
$ws.Range("F27").Value = 1701
$ws.Range("F31").Value = 668
$ws.Range("F32").Value = 886
$ws.Range("F34").Value = 2042
$ws.Range("F35").Value = 927
$ws.Range("F36").Value = 2066
$ws.Range("F38").Value = 511
$ws.Range("F39").Value = 232
$ws.Range("F40").Value = 857
$ws.Range("F42").Value = 966
$ws.Range("F43").Value = 831
$ws.Range("F44").Value = 1087
$ws.Range("F45").Value = 205
$ws.Range("F46").Value = 460
$ws.Range("F47").Value = 314
$ws.Range("F48").Value = 247

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 4
$ws.Range("F13").Value = 838
$ws.Range("F14").Value = 27
$ws.Range("F22").Value = 21

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 518
$ws.Range("F4").Value = 217
$ws.Range("F5").Value = 103
$ws.Range("F7").Value = 345
$ws.Range("F8").Value = 1382
$ws.Range("F12").Value = 212
$ws.Range("F13").Value = 150
$ws.Range("F15").Value = 1718
$ws.Range("F16").Value = 288
$ws.Range("F17").Value = 339
$ws.Range("F18").Value = 3639
$ws.Range("F19").Value = 35
$ws.Range("F20").Value = 4
$ws.Range("F24").Value = 1231
$ws.Range("F25").Value = 2881
$ws.Range("F27").Value = 1701
$ws.Range("F32").Value = 838
$ws.Range("F33").Value = 27
$ws.Range("F35").Value = 886
$ws.Range("F37").Value = 927
$ws.Range("F38").Value = 2066
$ws.Range("F39").Value = 511
$ws.Range("F40").Value = 232
$ws.Range("F41").Value = 857
$ws.Range("F42").Value = 966
$ws.Range("F43").Value = 831
$ws.Range("F44").Value = 1087
$ws.Range("F45").Value = 460
$ws.Range("F47").Value = 21
$ws.Range("F48").Value = 247
